$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "定保" (routine-maintenance) legend entry
$ws.Range("D2").ClearContents()

# Merge C2:D2 so the "維修" legend swatch now spans the freed-up width
$merged = $ws.Range("C2:D2")
$merged.Merge()

# Reflect the merged range as the active selection
$merged.Select()
